# edit.ps1
# Applies the "edited header banner and picture" commit:
#   1. Bumps the dark overlay rectangle ("header banner", shape id 6 on
#      slide 1) from 20% opacity to 40% opacity (alpha 20000 -> 40000).
#   2. Refreshes the cached text of the auto-update "date" placeholder
#      (datetimeFigureOut field) on every slide layout from 2018/10/31
#      to 2019/12/3, mirroring PowerPoint's own field re-cache on save.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Header banner rectangle: raise fill opacity from 20% to 40%
#    (PowerPoint's Fill.Transparency is the inverse of opacity, so
#    opacity 40% == transparency 0.6).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 6) {
        $shp.Fill.Transparency = 0.6
    }
}

# ---------------------------------------------------------------------
# 2) Date placeholder on every slide layout: update the cached text of
#    the datetimeFigureOut field to the new save date.
# ---------------------------------------------------------------------
$newDate = "2019/12/3"
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.HasTextFrame) {
            $isDatePlaceholder = $false
            if ($shape.Type -eq 14) {
                try {
                    if ($shape.PlaceholderFormat.Type -eq 16) {
                        $isDatePlaceholder = $true
                    }
                } catch {
                    $isDatePlaceholder = $false
                }
            }
            if ($isDatePlaceholder) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
